$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_val data (regenerated to filter save games)
$ws.Range("B2").Value = 0.04271373187048222
$ws.Range("C2").Value = 0.0000005461030343489881
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 1.28969408180815

$ws.Range("B3").Value = 1.455362044514542
$ws.Range("C3").Value = 10.34677158129881
$ws.Range("D3").Value = 261.3203778131603
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 283.3149644459102

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.7527432677738641
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 15.88780690183548

$ws.Range("B5").Value = 1.455362044514542
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 1133.036916526867
$ws.Range("G5").Value = 1136.900799921416
